# Arbeitszeit.xlsx — update the time-tracking entry for row 22.
#
# The diff shows:
#   * B22 changes from 3 to 4 (an extra hour logged)
#   * E2 (=SUM(B:B)) recalculates from 53.48 to 54.48 as a consequence
#   * the sheet's view/selection moves from D10 to B22, scrolled so row 9
#     is at the top of the visible window

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
[void]$ws.Activate()

# Bump the logged hours for the entry dated 2023-01-04 (row 22) from 3 to 4.
# E2's SUM(B:B) formula recalculates automatically afterwards.
$ws.Range("B22").Value = 4

# Match the author's final view state: scrolled down with B22 selected.
[void]$ws.Range("B22").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
